$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time, End Time, Run Time ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = '22:47:09 03-22-2017 Wednesday Pacific Daylight Time'
$wsVersion.Range("C11").Value = '22:47:44 03-22-2017 Wednesday Pacific Daylight Time'
$wsVersion.Range("C12").Value = '0:0:35 0 days'

# --- Sheet "3_3": update Deep Sleep (row 3) and BT Idle (row 4) stats ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.59
$ws33.Range("D3").Value = 0.612
$ws33.Range("E3").Value = 0.561
$ws33.Range("F3").Value = 0.014
$ws33.Range("H3").Value = '0.565942,0.575229,0.600093,0.602512,0.56993,0.593524,0.602056,0.601866,0.5769,0.607253,0.596302,0.592301,0.561017,0.597872,0.601303,0.600562,0.596751,0.576703,0.607253,0.608184,0.565942,0.570772,0.600195,0.601486,0.601037,0.571513,0.599726,0.603442,0.58116,0.597416,0.592118,0.59974,0.592953,0.572437,0.595175,0.596486,0.58099,0.571241,0.597505,0.604101,0.570025,0.573449,0.608279,0.601316,0.565841,0.598524,0.598897,0.570861,0.601214,0.598626,0.599359,0.570677,0.596669,0.603171,0.603259,0.564176,0.599828,0.596479,0.561771,0.568911,0.598157,0.603809,0.574013,0.590161,0.600487,0.611703,0.56724,0.601493,0.600188,0.600854,0.578381,0.597878,0.576431,0.588959,0.597967,0.579033,0.600385,0.597505,0.570025,0.605677,0.600664,0.564271,0.564455,0.602328,0.594631,0.604373,0.569468,0.605386,0.602518,0.570494,0.603531,0.572077,0.59546,0.603164,0.596391,0.589054,0.564917,0.602695,0.590446,0.572627'

$ws33.Range("C4").Value = 9.824
$ws33.Range("D4").Value = 9.841
$ws33.Range("E4").Value = 9.805
$ws33.Range("F4").Value = 0.009
$ws33.Range("H4").Value = '9.825227,9.812129,9.822903,9.828664,9.835064,9.8177,9.814826,9.832095,9.839242,9.826429,9.80685,9.824758,9.837856,9.827088,9.81789,9.812604,9.830987,9.835064,9.831259,9.820675,9.815478,9.837665,9.834126,9.825036,9.813243,9.830057,9.827822,9.830702,9.817516,9.817238,9.824663,9.839058,9.823922,9.805166,9.825689,9.832557,9.840546,9.816592,9.814085,9.815851,9.837292,9.828556,9.812414,9.813066,9.825967,9.838216,9.824296,9.807034,9.818725,9.83933,9.827088,9.821979,9.8088,9.819384,9.832095,9.824486,9.810743,9.81399,9.832006,9.828195,9.824574,9.816123,9.819384,9.836456,9.825036,9.824663,9.81668,9.824296,9.838589,9.828936,9.806925,9.819377,9.829962,9.829506,9.825505,9.81789,9.818542,9.833481,9.828005,9.81575,9.821137,9.824574,9.829316,9.821599,9.811857,9.821137,9.823631,9.833107,9.828474,9.811116,9.817896,9.831905,9.831639,9.813712,9.819187,9.826069,9.838766,9.82522,9.816871,9.814832'

# --- Sheet "1_8": update Deep Sleep (row 3) and BT Idle (row 4) stats ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0.001
$ws18.Range("D3").Value = 0.006
$ws18.Range("H3").Value = '-0.000762,-0.000263,-0.00802,0.001445,-0.006879,-0.001903,-0.000824,-0.004123,-0.002181,-0.003278,0.004201,0.000878,0.00253,-0.001632,-0.000263,0.000896,-0.003568,-0.000799,-0.003284,-0.001089,0.002567,-0.002489,-0.000509,-0.009117,-0.004955,-0.004413,0.003363,-0.004702,0.000866,-0.001342,0.001137,-0.001077,-0.002476,-0.002452,-0.003025,0.001982,0.002524,0.002278,-0.000793,-0.002458,-0.001656,-0.002452,0.003067,0.000594,-0.001928,-0.011078,-0.000799,0.003659,-0.005523,0.002512,-0.002168,-0.004388,-0.001354,0.000304,-0.004419,-0.002464,-0.000793,-0.004382,-0.004123,-0.001595,0.001969,-0.004123,-0.002187,-0.000812,0.000582,0.002555,-0.00773,-0.000522,-0.000516,-0.00276,-0.004123,0.00282,-0.000232,-0.005529,-0.005214,-0.00247,5.2e-05,0.001692,-0.001626,0.00446,0.001402,0.005564,0.003665,0.001427,0.001686,-0.002174,-0.000503,0.001988,-0.002754,-0.000516,-0.00215,3.9e-05,-0.002742,-0.010234,0.001698,-0.007194,-0.003019,-0.000781,0.001427,0.001427'

$ws18.Range("E4").Value = -0.009
$ws18.Range("H4").Value = '-0.000522,-0.003284,0.003918,-0.00247,-0.003845,0.001704,-0.003278,-0.001626,0.001155,-0.002476,0.005601,-0.005233,-0.005239,-0.001089,0.000304,-0.004111,0.008074,-0.002723,-0.002181,0.000576,-0.004098,2.1e-05,-0.002181,-0.005214,5.2e-05,3.9e-05,-0.008852,-0.004986,-0.003833,0.004485,-0.002729,-0.001626,0.002561,-0.000805,-0.00022,-0.003297,-0.005794,-0.003593,-0.001626,0.004201,-0.004653,0.000304,-0.00387,0.002808,-0.003556,0.0006,-0.002211,0.004738,-0.004659,-0.006355,-0.001064,-0.002464,-0.0044,-0.001077,-0.00247,0.001692,0.002278,-0.002458,-0.001613,-0.003568,0.003079,-0.000509,-0.005529,0.001704,0.002808,0.000866,-0.000787,-0.001077,-0.001885,-0.001342,-0.00691,-0.00329,0.000311,-0.000232,-0.004955,0.004763,-0.005541,-0.008563,0.001131,-0.003882,-0.001613,-0.002748,0.005891,-0.004394,0.002537,-0.003019,0.00393,-0.001897,-0.004394,-0.000793,0.001704,-0.002446,-0.004123,-0.002193,0.001143,0.001433,0.00171,-0.001632,0.001963,0.00475'
